$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 101
$ws.Range("D3").Value = 101
$ws.Range("D4").Value = 101
$ws.Range("D5").Value = 101
$ws.Range("D6").Value = 110
$ws.Range("D7").Value = 99
$ws.Range("D8").Value = 99
$ws.Range("D9").Value = 99
$ws.Range("D10").Value = 99
$ws.Range("D11").Value = 99
$ws.Range("D12").Value = 99
$ws.Range("D13").Value = 99

$ws.Range("D14").Select()
